# Testes com novos dados
# Add a new row (row 2) with two text values, "4004.0" and "5000.0",
# mirroring the data entered in the source spreadsheet. The values are
# numeric-looking but must be stored as literal text, so we briefly force
# a Text number format before assigning them (otherwise Excel's normal
# value-entry coercion would turn them into the numbers 4004 and 5000),
# then clear the formatting back to the sheet's default so the new cells
# don't pick up any extra styling.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:B2").NumberFormat = "@"
$ws.Range("A2").Value = "4004.0"
$ws.Range("B2").Value = "5000.0"
$ws.Range("A2:B2").ClearFormats()
